# Applies the "Deploying to gh-pages" update to the rx-supply-indicator
# StructureDefinition workbook:
#   - bump Version 5.0.0 -> 6.0.0
#   - bump Date to the new publish timestamp
#   - fill in the Publisher value ("Alvearie Team")
#   - replace the duplicated "Contact" metadata row with a "Jurisdiction" row
#     and remove the now-redundant extra row so the table is 20 rows again
#   - give the root Extension element its real Short/Definition text on the
#     Elements sheet instead of the generic placeholders

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B3").Value = "6.0.0"
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"
$meta.Range("B9").Value = "Alvearie Team"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Row 11 was a duplicate "Contact" row; delete it and shift everything
# below it up, which also drops the table back down to 20 rows.
$meta.Rows.Item(11).Delete()

$elements = $wb.Worksheets.Item("Elements")
$elements.Range("K2").Value = "Rx Supply Indicator"
$elements.Range("L2").Value = "Indicates whether the drug claim is a medical supply (Y) or a drug (N)"
